$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current header row (row 1), pushing
# everything down by two rows. The original row 1 (with its bold,
# bordered, centered header style) ends up on row 3.
$ws.Rows.Item(1).Insert()
$ws.Rows.Item(1).Insert()

# New row 1: numeric column index markers 0-12.
for ($i = 0; $i -lt 13; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $i
}

# Move the old header formatting from row 3 onto the new row 1 before
# stripping it off of row 3 (row 3 becomes a plain data-style row, same
# as the rest of the table).
$ws.Range("A3:M3").Copy()
$ws.Range("A1:M1").PasteSpecial(-4122)
$ws.Range("A3:M3").Style = "Normal"
$ws.Cells.Item(3, 10).ClearContents()

# New row 2: a mostly-blank row labelling the Flange/Drive sub-columns.
$ws.Cells.Item(2, 3).Value = "Flange"
$ws.Cells.Item(2, 6).Value = "Drive"
